# LOB1042.xlsx update
# The "Docentes responsáveis:" value row (old row 13, holding the
# B/C value "230696 - Carlos José Todero Peixoto") is removed, shifting
# every row below it up by one. On top of that shift, several of the
# value cells (column B/C) end up holding content taken from a
# different field than before, while the column-A labels stay aligned
# with their original row position after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old row 13 entirely (was: B13/C13 = "230696 - Carlos
#    José Todero Peixoto", no A13 label). Excel shifts rows 14-24 (and
#    their row heights) up to become rows 13-23.
$ws.Rows(13).Delete()

# 2) Fix up the value cells (B/C) that now hold different text than a
#    plain shift would produce. Column A labels already land correctly
#    after the delete, so only B/C need touching.

$ws.Range("B10").Value = "230696 - Carlos José Todero Peixoto"
$ws.Range("C10").Value = "230696 - Carlos José Todero Peixoto"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

$ws.Range("B18").Value = "230696 - Carlos José Todero Peixoto"
$ws.Range("C18").Value = "230696 - Carlos José Todero Peixoto"

$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

$ws.Range("B20").Value = "NF 5,0."
$ws.Range("C20").Value = "NF 5,0."

$ws.Range("B21").Value = "(NF+RC)/2 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 5,0, onde RC é uma prova de recuperação a ser aplicada."
